# Collapse the "valor total / parcelas" paragraph down to the new
# {{remuneracao}}{{parcelas}} ... {{percentual_ganho}} wording, removing the
# old installment-breakdown runs.

$d = $word.ActiveDocument

# Locate the target paragraph robustly (by distinctive placeholder text)
# rather than hard-coding a paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*{{valor_total}}*" -and $t -like "*{{numero_parcelas}}*" -and $t -like "*divididos em*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "target paragraph not found"
}

$pr = $target.Range
# Exclude the trailing paragraph mark from the replacement range.
$contentStart = $pr.Start
$contentEnd = $pr.End - 1
$replaceRange = $d.Range($contentStart, $contentEnd)

$newXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
<w:rtl w:val="0"/>
</w:rPr>
<w:t xml:space="preserve">{{remuneracao}}{{parcelas}}</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/>
<w:b w:val="0"/>
<w:i w:val="0"/>
<w:smallCaps w:val="0"/>
<w:strike w:val="0"/>
<w:color w:val="000000"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
<w:u w:val="none"/>
<w:vertAlign w:val="baseline"/>
<w:rtl w:val="0"/>
</w:rPr>
<w:t xml:space="preserve">. </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/>
<w:b w:val="1"/>
<w:i w:val="0"/>
<w:smallCaps w:val="0"/>
<w:strike w:val="0"/>
<w:color w:val="000000"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
<w:u w:val="none"/>
<w:vertAlign w:val="baseline"/>
<w:rtl w:val="0"/>
</w:rPr>
<w:t xml:space="preserve">{{percentual_ganho}}</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Calibri" w:cs="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri"/>
<w:b w:val="1"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
<w:rtl w:val="0"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$replaceRange.InsertXML($newXml)
